$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5346
$ws1.Range("F3").Value = 582
$ws1.Range("F4").Value = 11194
$ws1.Range("G4").Value = 62
$ws1.Range("F5").Value = 277
$ws1.Range("F6").Value = 584
$ws1.Range("F8").Value = 235
$ws1.Range("F9").Value = 962

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5346
$ws4.Range("F5").Value = 582
$ws4.Range("F7").Value = 11194
$ws4.Range("G7").Value = 62
$ws4.Range("F8").Value = 277
$ws4.Range("F9").Value = 584
$ws4.Range("F13").Value = 235
$ws4.Range("F14").Value = 962
